$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = 16975
$ws.Range("B41").Value = 12552
$ws.Range("B42").Value = 9928
$ws.Range("B43").Value = 7577
$ws.Range("B44").Value = 13563
$ws.Range("B45").Value = 17254
$ws.Range("B46").Value = 15730
